$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.855.48"
$ws.Range("E2").Value = "  +2.52%  "
$ws.Range("D3").Value = "3.320.59"
$ws.Range("E3").Value = "  +0.18%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'582.51"
$ws.Range("E5").Value = "  +3.37%  "
$ws.Range("D6").Value = "'182.86"
$ws.Range("E6").Value = "  -1.61%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "'0.590"
$ws.Range("E8").Value = "  +2.89%  "
$ws.Range("D9").Value = "3.316.91"
$ws.Range("E9").Value = "  +0.25%  "
$ws.Range("D10").Value = "'0.178"
$ws.Range("E10").Value = "  +1.33%  "
$ws.Range("D11").Value = "'0.580"
$ws.Range("E11").Value = "  +1.16%  "
$ws.Range("D12").Value = "'46.49"
$ws.Range("E12").Value = "  +0.92%  "
$ws.Range("E13").Value = "  +3.32%  "
$ws.Range("D14").Value = "'633.56"
$ws.Range("E14").Value = "  +6.91%  "
$ws.Range("D15").Value = "3.857.87"
$ws.Range("E15").Value = "  +0.32%  "
$ws.Range("E16").Value = "  +0.26%  "
$ws.Range("D17").Value = "67.982.53"
$ws.Range("E17").Value = "  +2.99%  "
$ws.Range("E18").Value = "  +1.53%  "
$ws.Range("D19").Value = "3.323.27"
$ws.Range("E19").Value = "  +0.29%  "
$ws.Range("D20").Value = "'17.71"
$ws.Range("E20").Value = "  +0.24%  "
$ws.Range("D21").Value = "'10.94"
$ws.Range("E21").Value = "  +0.23%  "
$ws.Range("E22").Value = "  +0.98%  "
$ws.Range("D23").Value = "'17.64"
$ws.Range("E23").Value = "  -2.29%  "
$ws.Range("D24").Value = "'5.05"
$ws.Range("E24").Value = "  +0.73%  "
$ws.Range("D25").Value = "'97.03"
$ws.Range("E25").Value = "  -1.18%  "
$ws.Range("E26").Value = "  +0.51%  "
$ws.Range("E27").Value = "  +2.70%  "
$ws.Range("D28").Value = "'9.59"
$ws.Range("E28").Value = "  +2.10%  "
$ws.Range("D29").Value = "'32.53"
$ws.Range("E29").Value = "  +6.30%  "
$ws.Range("D30").Value = "'8.61"
$ws.Range("E30").Value = "  +2.08%  "
$ws.Range("D31").Value = "'6.75"
$ws.Range("E31").Value = "  +1.77%  "
$ws.Range("D32").Value = "'593.60"
$ws.Range("E32").Value = "  +5.79%  "
$ws.Range("D33").Value = "3.946.97"
$ws.Range("E33").Value = "  +4.28%  "
$ws.Range("D34").Value = "'10.97"
$ws.Range("E34").Value = "  +1.34%  "
$ws.Range("E35").Value = "  +1.75%  "
$ws.Range("D36").Value = "'3.52"
$ws.Range("E36").Value = "  -4.76%  "
$ws.Range("E37").Value = "  -0.16%  "
$ws.Range("D38").Value = "'55.69"
$ws.Range("E38").Value = "  -0.39%  "
$ws.Range("E39").Value = "  +1.54%  "
$ws.Range("D40").Value = "'3.26"
$ws.Range("E40").Value = "  +3.80%  "
$ws.Range("E41").Value = "  +4.24%  "
$ws.Range("D42").Value = "'32.67"
$ws.Range("E42").Value = "  -1.83%  "
$ws.Range("E43").Value = "  +1.01%  "
$ws.Range("E44").Value = "  -0.01%  "
$ws.Range("D45").Value = "'0.339"
$ws.Range("E45").Value = "  +1.88%  "
$ws.Range("E46").Value = "  +0.96%  "
$ws.Range("E47").Value = "  +1.88%  "
$ws.Range("D48").Value = "'1.01"
$ws.Range("E48").Value = "  +0.67%  "
$ws.Range("D49").Value = "'1.39"
$ws.Range("E49").Value = "  +12.90%  "
$ws.Range("D50").Value = "'2.56"
$ws.Range("E50").Value = "  +1.25%  "
$ws.Range("D51").Value = "'130.95"
$ws.Range("E51").Value = "  +2.07%  "
